$wb = $excel.ActiveWorkbook

$wsTasks = $wb.Worksheets.Item(1)   # "Задачи"
$wsBacklog = $wb.Worksheets.Item(2) # "Бэклог задач"

# --- Sheet1 "Задачи": add finish date to row 8 (task "Сделать признак аналогичности серверному.") ---
$wsTasks.Range("D8").Value = 41975.586805555555
$wsTasks.Range("C8").Copy()
$wsTasks.Range("D8").PasteSpecial(-4122)

# --- Sheet2 "Бэклог задач": append three new backlog items ---
$wsBacklog.Range("B16").Value = "Убрать проверку switchOff, всё равно не используется"
$wsBacklog.Range("C16").Value = 41975.678472222222
$wsBacklog.Range("C15").Copy()
$wsBacklog.Range("C16").PasteSpecial(-4122)

$wsBacklog.Range("B17").Value = "Сделать внешние методы, либо in\out, либо setup, switchOn, switchOff"
$wsBacklog.Range("C17").Value = 41975.691666666666
$wsBacklog.Range("C15").Copy()
$wsBacklog.Range("C17").PasteSpecial(-4122)

$wsBacklog.Range("B18").Value = "Сделать валидацию in{Name} = out{Name}"
$wsBacklog.Range("C18").Value = 41975.691666666666
$wsBacklog.Range("C15").Copy()
$wsBacklog.Range("C18").PasteSpecial(-4122)

# --- Switch active sheet/selection from "Задачи" to "Бэклог задач" ---
$wsTasks.Range("B8").Select()
$wsBacklog.Activate()
$wsBacklog.Range("C18").Select()
